$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(8).Delete()
$ws.Cells.Item(2, 2).Value = 1.03
$ws.Cells.Item(2, 3).Value = 1.03
$ws.Cells.Item(2, 4).Value = 1.03
$ws.Cells.Item(3, 2).Value = 1.03
$ws.Cells.Item(3, 3).Value = 1.03
$ws.Cells.Item(3, 4).Value = 1.03
$ws.Cells.Item(3, 5).Value = 359.220000000004
$ws.Cells.Item(3, 6).Value = 82.02000000000092
$ws.Cells.Item(3, 7).Value = 323.6200000000035
$ws.Cells.Item(3, 9).Value = 41.01000000000046
$ws.Cells.Item(4, 2).Value = 1.03
$ws.Cells.Item(4, 3).Value = 1.03
$ws.Cells.Item(4, 4).Value = 1.03
$ws.Cells.Item(4, 5).Value = 844.4900000000283
$ws.Cells.Item(4, 8).Value = 460.190000000016
$ws.Cells.Item(4, 9).Value = 322.6400000000148
$ws.Cells.Item(5, 2).Value = 1.03
$ws.Cells.Item(5, 3).Value = 1.03
$ws.Cells.Item(5, 4).Value = 1.03
$ws.Cells.Item(5, 9).Value = 530.9900000000332
$ws.Cells.Item(6, 2).Value = 1.03
$ws.Cells.Item(6, 3).Value = 1.03
$ws.Cells.Item(6, 4).Value = 1.03
$ws.Cells.Item(7, 2).Value = 1.03
$ws.Cells.Item(7, 3).Value = 1.03
$ws.Cells.Item(7, 4).Value = 1.03
$ws.Cells.Item(7, 5).Value = 1408.090000000052
$ws.Cells.Item(7, 9).Value = 570.9900000000436
$ws.Cells.Item(7, 10).Value = 613.7100000000355

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(8).Delete()
$ws.Cells.Item(2, 2).Value = 1.03
$ws.Cells.Item(2, 3).Value = 1.03
$ws.Cells.Item(2, 4).Value = 1.03
$ws.Cells.Item(2, 7).Value = 323.6199999999999
$ws.Cells.Item(2, 10).Value = 161.8100000000004
$ws.Cells.Item(3, 2).Value = 1.03
$ws.Cells.Item(3, 3).Value = 1.03
$ws.Cells.Item(3, 4).Value = 1.03
$ws.Cells.Item(3, 5).Value = 844.4900000000181
$ws.Cells.Item(3, 7).Value = 697.0800000000148
$ws.Cells.Item(3, 9).Value = 322.640000000011
$ws.Cells.Item(3, 10).Value = 369.6700000000083
$ws.Cells.Item(4, 2).Value = 1.029999899221413
$ws.Cells.Item(4, 3).Value = 1.029999976988996
$ws.Cells.Item(4, 4).Value = 1.02999990920725
$ws.Cells.Item(4, 5).Value = 1266.520000000043
$ws.Cells.Item(4, 9).Value = 530.990000000025
$ws.Cells.Item(5, 2).Value = 1.029999755978219
$ws.Cells.Item(5, 3).Value = 1.029999840703543
$ws.Cells.Item(5, 4).Value = 1.029999800479636
$ws.Cells.Item(6, 2).Value = 1.029999636582067
$ws.Cells.Item(6, 3).Value = 1.029999737899343
$ws.Cells.Item(6, 4).Value = 1.029999697598747
$ws.Cells.Item(6, 5).Value = 1408.090000000033
$ws.Cells.Item(6, 9).Value = 570.9900000000358
$ws.Cells.Item(6, 10).Value = 613.7100000000273
$ws.Cells.Item(7, 2).Value = 1.029999596865019
$ws.Cells.Item(7, 3).Value = 1.029999715455052
$ws.Cells.Item(7, 4).Value = 1.029999669543382

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(8).Delete()
$ws.Cells.Item(2, 2).Value = 1.03
$ws.Cells.Item(2, 3).Value = 1.03
$ws.Cells.Item(2, 4).Value = 1.03
$ws.Cells.Item(3, 2).Value = 1.022468000387581
$ws.Cells.Item(3, 3).Value = 1.032998423494961
$ws.Cells.Item(3, 4).Value = 1.02035317770982
$ws.Cells.Item(3, 5).Value = 625.8800000000097
$ws.Cells.Item(4, 2).Value = 1.006558041103589
$ws.Cells.Item(4, 3).Value = 1.023124404042142
$ws.Cells.Item(4, 4).Value = 1.012346853434017
$ws.Cells.Item(4, 6).Value = 593.420000000008
$ws.Cells.Item(4, 7).Value = 629.3600000000085
$ws.Cells.Item(4, 8).Value = 412.3400000000054
$ws.Cells.Item(4, 10).Value = 326.8200000000046
$ws.Cells.Item(5, 2).Value = 0.9943485024708942
$ws.Cells.Item(5, 3).Value = 1.016512310413684
$ws.Cells.Item(5, 4).Value = 1.002920343872736
$ws.Cells.Item(5, 6).Value = 593.4200000000278
$ws.Cells.Item(5, 7).Value = 629.36000000003
$ws.Cells.Item(5, 8).Value = 412.340000000019
$ws.Cells.Item(5, 9).Value = 339.5700000000164
$ws.Cells.Item(5, 10).Value = 326.8200000000161
$ws.Cells.Item(6, 2).Value = 0.9899730981916601
$ws.Cells.Item(6, 3).Value = 1.016129781667785
$ws.Cells.Item(6, 4).Value = 1.00000771969154
$ws.Cells.Item(6, 6).Value = 593.420000000008
$ws.Cells.Item(6, 7).Value = 629.3600000000085
$ws.Cells.Item(6, 8).Value = 412.3400000000054
$ws.Cells.Item(6, 10).Value = 326.8200000000046
$ws.Cells.Item(7, 2).Value = 0.9899729739680746
$ws.Cells.Item(7, 3).Value = 1.016129677460319
$ws.Cells.Item(7, 4).Value = 1.000007613725063
$ws.Cells.Item(7, 6).Value = 593.420000000008
$ws.Cells.Item(7, 7).Value = 629.3600000000085
$ws.Cells.Item(7, 8).Value = 412.3400000000054
$ws.Cells.Item(7, 10).Value = 326.8200000000046

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(8).Delete()
$ws.Cells.Item(2, 2).Value = 1.03
$ws.Cells.Item(2, 3).Value = 1.03
$ws.Cells.Item(2, 4).Value = 1.03
$ws.Cells.Item(3, 2).Value = 1.024508217831734
$ws.Cells.Item(3, 3).Value = 1.023915208218881
$ws.Cells.Item(3, 4).Value = 1.02431898186023
$ws.Cells.Item(4, 2).Value = 0.9998872004308923
$ws.Cells.Item(4, 3).Value = 1.014203597998166
$ws.Cells.Item(4, 4).Value = 1.001165073778819
$ws.Cells.Item(5, 2).Value = 0.9772571126707913
$ws.Cells.Item(5, 3).Value = 1.002693840277917
$ws.Cells.Item(5, 4).Value = 0.9886917718568555
$ws.Cells.Item(6, 2).Value = 0.9646767948496874
$ws.Cells.Item(6, 3).Value = 0.9959461172208336
$ws.Cells.Item(6, 4).Value = 0.9790375295327033
$ws.Cells.Item(7, 2).Value = 0.9601661883749877
$ws.Cells.Item(7, 3).Value = 0.9955556862259944
$ws.Cells.Item(7, 4).Value = 0.9760536396747709
$ws.Cells.Item(7, 5).Value = 478.0899354777617
$ws.Cells.Item(7, 6).Value = 372.9789915816977
$ws.Cells.Item(7, 7).Value = 395.6187394771221
$ws.Cells.Item(7, 8).Value = 244.357987801907
$ws.Cells.Item(7, 9).Value = 224.1894957908488
$ws.Cells.Item(7, 10).Value = 213.8393697385611

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Rows.Item(8).Delete()
$ws.Cells.Item(2, 2).Value = 1.03
$ws.Cells.Item(2, 3).Value = 1.03
$ws.Cells.Item(2, 4).Value = 1.03
$ws.Cells.Item(3, 2).Value = 1.02712058454737
$ws.Cells.Item(3, 3).Value = 1.027216568872743
$ws.Cells.Item(3, 4).Value = 1.026810225653305
$ws.Cells.Item(4, 2).Value = 1.019779125975322
$ws.Cells.Item(4, 3).Value = 1.020521251826379
$ws.Cells.Item(4, 4).Value = 1.019611249557229
$ws.Cells.Item(5, 2).Value = 0.9950411001339402
$ws.Cells.Item(5, 3).Value = 1.010777032782102
$ws.Cells.Item(5, 4).Value = 0.9963479455386173
$ws.Cells.Item(6, 2).Value = 0.9722982266958458
$ws.Cells.Item(6, 3).Value = 0.9992278064278642
$ws.Cells.Item(6, 4).Value = 0.9838135714463
$ws.Cells.Item(7, 2).Value = 0.9596529038624853
$ws.Cells.Item(7, 3).Value = 0.992456518047324
$ws.Cells.Item(7, 4).Value = 0.9741109833371012
$ws.Cells.Item(7, 5).Value = 101.57
$ws.Cells.Item(7, 8).Value = 50.78

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Rows.Item(8).Delete()
$ws.Cells.Item(2, 2).Value = 1.03
$ws.Cells.Item(2, 3).Value = 1.03
$ws.Cells.Item(2, 4).Value = 1.03
$ws.Cells.Item(3, 2).Value = 1.024508217831734
$ws.Cells.Item(3, 3).Value = 1.023915208218881
$ws.Cells.Item(3, 4).Value = 1.02431898186023
$ws.Cells.Item(4, 2).Value = 0.9998872127478373
$ws.Cells.Item(4, 3).Value = 1.014203608453777
$ws.Cells.Item(4, 4).Value = 1.001165084378498
$ws.Cells.Item(5, 2).Value = 0.9772571294650578
$ws.Cells.Item(5, 3).Value = 1.002693853162426
$ws.Cells.Item(5, 4).Value = 0.9886917855172327
$ws.Cells.Item(6, 2).Value = 0.9646768118629674
$ws.Cells.Item(6, 3).Value = 0.9959461301926376
$ws.Cells.Item(6, 4).Value = 0.9790375433277848
$ws.Cells.Item(7, 2).Value = 0.9601662054681915
$ws.Cells.Item(7, 3).Value = 0.9955556992028856
$ws.Cells.Item(7, 4).Value = 0.9760536535120252

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Rows.Item(8).Delete()
$ws.Cells.Item(2, 2).Value = 1.03
$ws.Cells.Item(2, 3).Value = 1.03
$ws.Cells.Item(2, 4).Value = 1.03
$ws.Cells.Item(2, 5).Value = 281.4199999999997
$ws.Cells.Item(3, 2).Value = 1.018936692591816
$ws.Cells.Item(3, 3).Value = 1.034124074315433
$ws.Cells.Item(3, 4).Value = 1.016800756829455
$ws.Cells.Item(3, 5).Value = 281.4199999999997
$ws.Cells.Item(4, 2).Value = 0.9978297910078687
$ws.Cells.Item(4, 3).Value = 1.018872909387555
$ws.Cells.Item(4, 4).Value = 1.004927162186314
$ws.Cells.Item(4, 5).Value = 281.4199999999997
$ws.Cells.Item(5, 2).Value = 0.9855121362569483
$ws.Cells.Item(5, 3).Value = 1.012233055065915
$ws.Cells.Item(5, 4).Value = 0.9954304029011509
$ws.Cells.Item(5, 5).Value = 281.4199999999997
$ws.Cells.Item(6, 2).Value = 0.9810973293615369
$ws.Cells.Item(6, 3).Value = 1.011848910839918
$ws.Cells.Item(6, 4).Value = 0.9924958015511011
$ws.Cells.Item(6, 5).Value = 281.4199999999997
$ws.Cells.Item(7, 2).Value = 0.9810972040141283
$ws.Cells.Item(7, 3).Value = 1.011848806191577
$ws.Cells.Item(7, 4).Value = 0.9924956947825938
$ws.Cells.Item(7, 5).Value = 281.4199999999997
